$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M9").ClearContents()
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0

$ws.Range("H43").Value = 1764.2858
$ws.Range("I43").Value = 1977.7778
$ws.Range("K43").Value = 1977.7778
$ws.Range("M43").Value = -1908.7778

$ws.Range("H74").Value = 9506.1875
$ws.Range("I74").Value = 10309.9
$ws.Range("K74").Value = 10309.9
$ws.Range("M74").Value = -9373.9

$ws.Range("H77").Value = 9506.1875
$ws.Range("I77").Value = 10309.9
$ws.Range("K77").Value = 51549.5
$ws.Range("M77").Value = -46869.5

$ws.Range("H92").Value = 345.30768
$ws.Range("I92").Value = 320.85715
$ws.Range("J92").Value = 373.83334
$ws.Range("K92").Value = 320.85715
$ws.Range("L92").Value = 373.83334
$ws.Range("M92").Value = 927.14285
$ws.Range("N92").Value = -2869.83334

$ws.Range("H103").Value = 448.6
$ws.Range("J103").Value = 572
$ws.Range("L103").Value = 1716
$ws.Range("N103").Value = -2888

$ws.Range("H106").Value = 3544.5
$ws.Range("J106").Value = 4995
$ws.Range("L106").Value = 4995
$ws.Range("N106").Value = -6257

$ws.Range("H107").Value = 33986.234
$ws.Range("I107").Value = 36363.855
$ws.Range("J107").Value = 699.5
$ws.Range("K107").Value = 36363.855
$ws.Range("L107").Value = 699.5
$ws.Range("M107").Value = -34443.855
$ws.Range("N107").Value = -4539.5

$ws.Range("H115").Value = 417.25
$ws.Range("I115").Value = 417.25
$ws.Range("K115").Value = 1251.75
$ws.Range("M115").Value = 315.25

$ws.Range("H128").Value = 88703
$ws.Range("J128").Value = 88703
$ws.Range("L128").Value = 88703
$ws.Range("N128").Value = -98663

$ws.Range("H132").Value = 14768.0625
$ws.Range("I132").Value = 2306.3572
$ws.Range("K132").Value = 6919.071599999999
$ws.Range("M132").Value = -4389.071599999999

$ws.Range("H137").Value = 3772.9429
$ws.Range("J137").Value = 3099.3333
$ws.Range("L137").Value = 9297.999899999999
$ws.Range("N137").Value = -14397.9999

$ws.Range("H141").Value = 3093.1
$ws.Range("I141").Value = 2734.8125
$ws.Range("K141").Value = 8204.4375
$ws.Range("M141").Value = -3024.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 30124.5
$ws.Range("I6").Value = 30124.5
$ws.Range("K6").Value = 30124.5
$ws.Range("M6").Value = -29951.5

$ws.Range("H97").Value = 1144.7222
$ws.Range("I97").Value = 1227.6
$ws.Range("J97").Value = 730.3333
$ws.Range("K97").Value = 1227.6
$ws.Range("L97").Value = 730.3333
$ws.Range("M97").Value = -731.5999999999999
$ws.Range("N97").Value = -1722.3333

$ws.Range("H106").Value = 50000
$ws.Range("J106").Value = 50000
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524

$ws.Range("H110").Value = 279336.56
$ws.Range("I110").Value = 279336.56
$ws.Range("K110").Value = 279336.56
$ws.Range("M110").Value = -277291.56

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I20").Value = 55558060
$ws.Range("J20").Value = 2335.4707
$ws.Range("K20").Value = 55558060
$ws.Range("L20").Value = 2335.4707
$ws.Range("M20").Value = -55557813
$ws.Range("N20").Value = -2829.4707

$ws.Range("H29").Value = 3016
$ws.Range("I29").Value = 3016
$ws.Range("K29").Value = 3016
$ws.Range("M29").Value = -2727

$ws.Range("H76").Value = 29157.25
$ws.Range("I76").Value = 6000
$ws.Range("K76").Value = 6000
$ws.Range("M76").Value = -5685

$ws.Range("H79").Value = 29157.25
$ws.Range("I79").Value = 6000
$ws.Range("K79").Value = 6000
$ws.Range("M79").Value = -4908

$ws.Range("H86").Value = 1001994.5
$ws.Range("I86").Value = 1215850.5
$ws.Range("J86").Value = 3999.6667
$ws.Range("K86").Value = 1215850.5
$ws.Range("L86").Value = 3999.6667
$ws.Range("M86").Value = -1214727.5
$ws.Range("N86").Value = -6245.6667

$ws.Range("H88").Value = 13999
$ws.Range("J88").Value = 13999
$ws.Range("L88").Value = 13999
$ws.Range("N88").Value = -14811

$ws.Range("H89").Value = 1001994.5
$ws.Range("I89").Value = 1215850.5
$ws.Range("J89").Value = 3999.6667
$ws.Range("K89").Value = 6079252.5
$ws.Range("L89").Value = 19998.3335
$ws.Range("M89").Value = -6073636.5
$ws.Range("N89").Value = -31230.3335

$ws.Range("H91").Value = 13999
$ws.Range("J91").Value = 13999
$ws.Range("L91").Value = 13999
$ws.Range("N91").Value = -16807

$ws.Range("H107").Value = 436871.22
$ws.Range("I107").Value = 1607.1177
$ws.Range("J107").Value = 1670119.5
$ws.Range("K107").Value = 1607.1177
$ws.Range("L107").Value = 1670119.5
$ws.Range("M107").Value = 312.8823
$ws.Range("N107").Value = -1673959.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 12000
$ws.Range("I56").Value = 12000
$ws.Range("K56").Value = 12000
$ws.Range("M56").Value = -11155

$ws.Range("H88").Value = 76999
$ws.Range("J88").Value = 76999
$ws.Range("L88").Value = 76999
$ws.Range("N88").Value = -77811

$ws.Range("H91").Value = 76999
$ws.Range("J91").Value = 76999
$ws.Range("L91").Value = 76999
$ws.Range("N91").Value = -79807

$ws.Range("H105").Value = 1010.1667
$ws.Range("J105").Value = 1010.5
$ws.Range("L105").Value = 1010.5
$ws.Range("N105").Value = -4504.5

$ws.Range("H107").Value = 485.57144
$ws.Range("I107").Value = 446.46155
$ws.Range("K107").Value = 446.46155
$ws.Range("M107").Value = 1473.53845

$ws.Range("H122").Value = 3547.1667
$ws.Range("I122").Value = 2454.6
$ws.Range("J122").Value = 4327.5713
$ws.Range("K122").Value = 7363.799999999999
$ws.Range("L122").Value = 12982.7139
$ws.Range("M122").Value = -4913.799999999999
$ws.Range("N122").Value = -17882.7139

$ws.Range("H132").Value = 2507.8
$ws.Range("I132").Value = 2250
$ws.Range("K132").Value = 6750
$ws.Range("M132").Value = -4220

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 254.53334
$ws.Range("I2").Value = 66.5
$ws.Range("K2").Value = 399
$ws.Range("M2").Value = -286

$ws.Range("H48").Value = 2997.5
$ws.Range("J48").Value = 2997.5
$ws.Range("L48").Value = 8992.5
$ws.Range("N48").Value = -9492.5

$ws.Range("H87").Value = 6795.8
$ws.Range("I87").Value = 6795.8
$ws.Range("K87").Value = 20387.4
$ws.Range("M87").Value = -19139.4

$ws.Range("H90").Value = 6795.8
$ws.Range("I90").Value = 6795.8
$ws.Range("K90").Value = 61162.2
$ws.Range("M90").Value = -54922.2

$ws.Range("H97").Value = 599.2
$ws.Range("J97").Value = 67.5
$ws.Range("L97").Value = 202.5
$ws.Range("N97").Value = -1194.5

$ws.Range("H129").Value = 125745.125
$ws.Range("I129").Value = 488
$ws.Range("J129").Value = 501516.5
$ws.Range("K129").Value = 1464
$ws.Range("L129").Value = 1504549.5
$ws.Range("M129").Value = 3536
$ws.Range("N129").Value = -1514549.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 49500
$ws.Range("J40").Value = 49500
$ws.Range("L40").Value = 49500
$ws.Range("N40").Value = -49802

$ws.Range("H70").Value = 8453.4
$ws.Range("I70").Value = 6407.6665
$ws.Range("J70").Value = 16636.334
$ws.Range("K70").Value = 6407.6665
$ws.Range("L70").Value = 16636.334
$ws.Range("M70").Value = -6137.6665
$ws.Range("N70").Value = -17176.334

$ws.Range("H73").Value = 8453.4
$ws.Range("I73").Value = 6407.6665
$ws.Range("J73").Value = 16636.334
$ws.Range("K73").Value = 6407.6665
$ws.Range("L73").Value = 16636.334
$ws.Range("M73").Value = -5471.6665
$ws.Range("N73").Value = -18508.334

$ws.Range("H80").Value = 1542662.9
$ws.Range("I80").Value = 1004461.2
$ws.Range("J80").Value = 3336668.8
$ws.Range("K80").Value = 1004461.2
$ws.Range("L80").Value = 3336668.8
$ws.Range("M80").Value = -1003463.2
$ws.Range("N80").Value = -3338664.8

$ws.Range("H83").Value = 1542662.9
$ws.Range("I83").Value = 1004461.2
$ws.Range("J83").Value = 3336668.8
$ws.Range("K83").Value = 5022306
$ws.Range("L83").Value = 16683344
$ws.Range("M83").Value = -5017314
$ws.Range("N83").Value = -16693328

$ws.Range("H102").Value = 2051.842
$ws.Range("I102").Value = 998.9286
$ws.Range("K102").Value = 998.9286
$ws.Range("M102").Value = 623.0714

$ws.Range("H113").Value = 1443972.1
$ws.Range("I113").Value = 2500823
$ws.Range("J113").Value = 34837.668
$ws.Range("K113").Value = 2500823
$ws.Range("L113").Value = 34837.668
$ws.Range("M113").Value = -2498653
$ws.Range("N113").Value = -39177.668

$ws.Range("H132").Value = 47040.96
$ws.Range("I132").Value = 7001.1904
$ws.Range("J132").Value = 257249.75
$ws.Range("K132").Value = 21003.5712
$ws.Range("L132").Value = 771749.25
$ws.Range("M132").Value = -18473.5712
$ws.Range("N132").Value = -776809.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2077.9167
$ws.Range("J46").Value = 662.6667
$ws.Range("L46").Value = 662.6667
$ws.Range("N46").Value = -1038.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 1000000
$ws.Range("J7").Value = 1000000
$ws.Range("L7").Value = 1000000
$ws.Range("N7").Value = -1000226

$ws.Range("H132").Value = 105703.3
$ws.Range("I132").Value = 4306.6
$ws.Range("K132").Value = 12919.8
$ws.Range("M132").Value = -10389.8

$ws.Range("H136").Value = 12423085
$ws.Range("I136").Value = 15876467
$ws.Range("J136").Value = 336248.84
$ws.Range("K136").Value = 47629401
$ws.Range("L136").Value = 1008746.52
$ws.Range("M136").Value = -47626851
$ws.Range("N136").Value = -1013846.52
